$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 7), continuing the daily expense log.
# Copy formatting from the row above (row 6) so the date cell keeps
# the existing date-number-format style instead of creating a new one.
$ws.Range("A6:M6").Copy($ws.Range("A7:M7"))

$ws.Range("A7").Value = 43795
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 6.5
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 3

$ws.Range("N7").Select()
